# Generate Report for Handoff
# Refresh the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# c7e9664d-cf34-4455-a087-7566a75630df.md file (row 7 on every sheet) to reflect a
# freshly generated handoff report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-26 02:42:07"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-26 02:41:58"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-26 02:42:07"
